# ---------------------------------------------------------------------------
# "oct 1 stuff - spreadsheet column name changes" (#2225)
#
# 1. Instructions sheet: rewrite the instructional copy, renumber the field
#    validation notes to match the new (human readable) header names, and
#    add a new "no field may be blank" note.
# 2. Forecast Report sheet: rename/reorder headers to the new human-readable
#    names, bold the header row, widen a couple of columns, zoom to 120%,
#    and attach dropdown (list) data validations sourced from a new sheet.
# 3. Add a new hidden "Dropdowns" sheet that backs those list validations.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Instructions sheet
# ---------------------------------------------------------------------------
$instructions = $wb.Worksheets.Item("Instructions")
$instructions.Unprotect()

$instructions.Range("A1").Value = 'Please fill out the "Forecast Report" sheet'
$instructions.Range("A3").Value = 'Please do not alter the name of the "Forecast Report" sheet'
$instructions.Range("A5").Value = 'Please do not alter any of the header names in the "Forecast Report" sheet'
$instructions.Range("A7").Value = 'Please note that no field in the "Forecast Report" sheet may be blank'

$instructions.Range("A9").Value = "Please note that:"
$instructions.Range("A10").Value = '(1) "Model Year" should be a 4 digit integer'
$instructions.Range("A11").Value = '(2) "Make" should be no more than 250 characters'
$instructions.Range("A12").Value = '(3) "Model" should be no more than 250 characters'
$instructions.Range("A13").Value = '(4) "Type" should be exactly one of: BEV, PHEV, FCEV, EREV'
$instructions.Range("A14").Value = '(5) "Range" should be a real number with no more than 2 decimal places'
$instructions.Range("A15").Value = '(6) "ZEV Class" should be a single, uppercase letter'
$instructions.Range("A16").Value = '(7) "Vehicle Class and Interior Volume" should be no more than 250 characters'
$instructions.Range("A17").Value = '(8) "Total ZEVs Supplied" should be an integer'

$instructions.Protect()

# ---------------------------------------------------------------------------
# 2. New hidden "Dropdowns" sheet (added after "Forecast Report")
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$dropdowns = $wb.Worksheets.Add([Type]::Missing, $lastSheet)
$dropdowns.Name = "Dropdowns"

$dropdowns.Range("A1").Value = "Type"
$dropdowns.Range("B1").Value = "ZEV Class"
$dropdowns.Range("C1").Value = "Vehicle Class and Interior Volume"
$dropdowns.Range("A1:C1").Font.Bold = $true

$dropdowns.Range("A2").Value = "BEV"
$dropdowns.Range("A3").Value = "PHEV"
$dropdowns.Range("A4").Value = "EREV"
$dropdowns.Range("A5").Value = "FCEV"

$dropdowns.Range("B2").Value = "A"
$dropdowns.Range("B3").Value = "B"
$dropdowns.Range("B4").Value = "C"

$dropdowns.Range("C2").Value = "Two-seater"
$dropdowns.Range("C3").Value = "Minicompact (less than 85 cu. ft.)"
$dropdowns.Range("C4").Value = "Subcompact (85–99 cu. ft.)"
$dropdowns.Range("C5").Value = "Compact (100–109 cu. ft.)"
$dropdowns.Range("C6").Value = "Mid-size (110–119 cu. ft.)"
$dropdowns.Range("C7").Value = "Full-size (120 cu. ft. or more)"
$dropdowns.Range("C8").Value = "Station wagon: Small (less than 130 cu. ft.)"
$dropdowns.Range("C9").Value = "Station wagon: Mid-size (130–159 cu. ft.)"
$dropdowns.Range("C10").Value = "Pickup truck: Small (less than 2,722 kg)"
$dropdowns.Range("C11").Value = "Pickup truck: Standard (2,722–3,856 kg)"
$dropdowns.Range("C12").Value = "Sport utility vehicle: Small (less than 2,722 kg)"
$dropdowns.Range("C13").Value = "Sport utility vehicle: Standard (2,722–4,536 kg)"
$dropdowns.Range("C14").Value = "Minivan (less than 3,856 kg)"
$dropdowns.Range("C15").Value = "Van: Cargo (less than 3,856 kg)"
$dropdowns.Range("C16").Value = "Van: Passenger (less than 4,536 kg)"
$dropdowns.Range("C17").Value = "Special purpose vehicle (less than 3,856 kg)"
$dropdowns.Range("C18").Value = "Other/TBD"

$dropdowns.Columns.Item(2).ColumnWidth = 8.166666666666666
$dropdowns.Columns.Item(3).ColumnWidth = 37.33333333333333

$dropdowns.Activate()
$dropdowns.Range("C7").Select()
$excel.ActiveWindow.Zoom = 130
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1

$dropdowns.Visible = 0

# ---------------------------------------------------------------------------
# 3. Forecast Report sheet - renamed / reordered headers + validations
# ---------------------------------------------------------------------------
$report = $wb.Worksheets.Item("Forecast Report")

$report.Range("A1").Value = "Model Year"
$report.Range("B1").Value = "Make"
$report.Range("C1").Value = "Model"
$report.Range("D1").Value = "Type"
$report.Range("E1").Value = "Range"
$report.Range("F1").Value = "ZEV Class"
$report.Range("G1").Value = "Vehicle Class and Interior Volume"
$report.Range("H1").Value = "Total ZEVs Supplied"

$report.Range("A1:H1").Font.Bold = $true

$report.Columns.Item(7).ColumnWidth = 28.498697916666668
$report.Columns.Item(8).ColumnWidth = 16.498697916666668

$report.Range("D2:D200").Validation.Add(3, 1, 1, "=Dropdowns!`$A`$2:`$A`$5")
$report.Range("F2:F200").Validation.Add(3, 1, 1, "=Dropdowns!`$B`$2:`$B`$4")
$report.Range("G2:G200").Validation.Add(3, 1, 1, "=Dropdowns!`$C`$2:`$C`$18")

$report.Activate()
$report.Range("A2").Select()
$excel.ActiveWindow.Zoom = 120

# ---------------------------------------------------------------------------
# Leave "Instructions" as the active / selected sheet, matching the source
# workbook's original tab selection.
# ---------------------------------------------------------------------------
$instructions.Activate()

Write-Output "Applied Oct 1 column-name changes."
